# Add image start/end dates from SD cards 83, 55, 84, 41 and 17
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Card 83 (row 7, Point_ID 1 / camera at row 7)
$ws.Range("G7").Value = 44074
$ws.Range("H7").Value = 44421

# Card 84 (row 8) - SD_card label corrected from "?" to "084?", plus dates
$ws.Range("D8").Value = "084?"
$ws.Range("F8").Value = 44421
$ws.Range("G8").Value = 42735
$ws.Range("H8").Value = 42735

# Card 55 (row 11)
$ws.Range("G11").Value = 44075
$ws.Range("H11").Value = 42737

# Card 41 (row 12)
$ws.Range("F12").Value = 44421
$ws.Range("G12").Value = 44076
$ws.Range("H12").Value = 42775

# Card 17 (row 14)
$ws.Range("F14").Value = 44421
$ws.Range("G14").Value = 44076
$ws.Range("H14").Value = 44083

# Page setup (paper size / orientation, matches the pageSetup element added to the sheet)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to match the final saved cursor position
$ws.Range("J11").Select()
